$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of data (Fecha = 44512) is being inserted into the table right
# before the existing row that used to be the first row of this block
# (row 136), pushing that row and everything after it down by two rows
# (one row per "Primera"/"Segunda" quality pair). The two rows that fall off
# the bottom of the previous range end up appended as the new last two rows
# (262/263).

# 1) Insert two blank rows at 136:137 - everything from row 136 down
#    (including the old last rows 260:261) shifts down to 138:139 ... 262:263.
$ws.Rows("136:137").Insert()

# 2) The new blank rows should start out as a copy of the record that just
#    got shifted down into 138:139 (same Mercado/Region/Categoria/Calidad/
#    Unidad/Origen/Clasificacion, etc.), then we only touch the two columns
#    that actually differ for the new week: Fecha (D) and Volumen (J).
$ws.Range("A138:R139").Copy()
$ws.Range("A136").PasteSpecial()

# 3) Overwrite Fecha/Volumen for the freshly inserted "Primera" (136) and
#    "Segunda" (137) rows with the new week's reported values.
$ws.Range("D136").Value = 44512
$ws.Range("J136").Value = 3000

$ws.Range("D137").Value = 44512
$ws.Range("J137").Value = 1400
